# Edit: change "A first publication" -> "The first publication"
# (Word's grammar-checker proofErr markers around "A first" become stale
#  once the text changes, so Find/Replace naturally removes them as the
#  run(s) get rewritten.)

$d = $word.ActiveDocument

# Make sure the edit is applied as a plain text change, not as a tracked
# revision (we want clean <w:r>/<w:t> runs, matching the target diff).
$d.TrackRevisions = $false

$d.Content.Find.Execute("A first publication", $true, $true, $false, $false,
                         $false, $true, 1, $false, "The first publication", 2)
